$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('word', 'image', 'category'),
    @('wenden', 'none', 'none'),
    @('schaden', 'flower/flower027.jpg', 'flower'),
    @('fordern', 'dog/dog028.jpg', 'dog'),
    @('rufen', 'none', 'none'),
    @('binden', 'flower/flower004.jpg', 'flower'),
    @('machen', 'dog/dog026.jpg', 'dog'),
    @('enden', 'none', 'none'),
    @('zählen', 'dog/dog001.jpg', 'dog'),
    @('werden', 'dog/dog009.jpg', 'dog'),
    @('opfern', 'none', 'none'),
    @('achten', 'dog/dog008.jpg', 'dog'),
    @('segeln', 'flower/flower008.jpg', 'flower'),
    @('schalten', 'none', 'none'),
    @('treten', 'flower/flower030.jpg', 'flower'),
    @('landen', 'dog/dog031.jpg', 'dog'),
    @('klagen', 'none', 'none'),
    @('schneiden', 'flower/flower031.jpg', 'flower'),
    @('wehtun', 'flower/flower003.jpg', 'flower'),
    @('sparen', 'none', 'none'),
    @('tanzen', 'flower/flower029.jpg', 'flower'),
    @('stopfen', 'dog/dog002.jpg', 'dog'),
    @('hören', 'none', 'none'),
    @('packen', 'dog/dog003.jpg', 'dog'),
    @('süßen', 'dog/dog000.jpg', 'dog'),
    @('drohen', 'none', 'none'),
    @('seufzen', 'dog/dog016.jpg', 'dog'),
    @('wundern', 'flower/flower005.jpg', 'flower'),
    @('orten', 'none', 'none'),
    @('retten', 'flower/flower014.jpg', 'flower'),
    @('bellen', 'dog/dog012.jpg', 'dog'),
    @('weigern', 'none', 'none'),
    @('füttern', 'dog/dog021.jpg', 'dog'),
    @('ruhen', 'flower/flower026.jpg', 'flower'),
    @('ächzen', 'none', 'none'),
    @('pflegen', 'dog/dog005.jpg', 'dog'),
    @('wüten', 'flower/flower019.jpg', 'flower'),
    @('kosten', 'none', 'none'),
    @('passen', 'flower/flower023.jpg', 'flower'),
    @('steuern', 'flower/flower006.jpg', 'flower'),
    @('dauern', 'none', 'none'),
    @('reizen', 'flower/flower000.jpg', 'flower'),
    @('ärgern', 'dog/dog006.jpg', 'dog'),
    @('stören', 'none', 'none'),
    @('wagen', 'dog/dog022.jpg', 'dog'),
    @('stillen', 'dog/dog015.jpg', 'dog'),
    @('bremsen', 'none', 'none'),
    @('zahlen', 'flower/flower015.jpg', 'flower'),
    @('triefen', 'flower/flower012.jpg', 'flower')
)

for ($r = 0; $r -lt $data.Count; $r++) {
    for ($c = 0; $c -lt 3; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}
